# plotEIC methods for fGroupsSet
#
# Inserts a new tracking row for "getEICsForFGroups" above "getFeatures"
# (pushing every row below it down by one), and marks the existing
# "plotEIC" row as done (column G) in addition to its current
# "implement" (column D) marker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 19 (shifts rows 19..53 down to 20..54)
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row: getEICsForFGroups
$ws.Range("A19").Value = "getEICsForFGroups"
$ws.Range("D19").Value = "X"
$ws.Range("F19").Value = "X"
$ws.Range("G19").Value = "X"

# plotEIC (now row 34) is also done
$ws.Range("G34").Value = "X"

# Update the saved selection to match where editing ended up
$ws.Range("G35").Select()
